# Apply recalibrated extrapolation values, excluding noisy sub-$5 price points.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value
$updates = @{
    2  = @{ D = 117592.8751483996;  E = -0.03587617850236063; F = 0.1677052604798103;  G = -1.456981444281443; H = 13.47495016257348 }
    4  = @{ D = 119217.0483713247;  E = -0.03680249953377605; F = 0.2459486504646433;  G = -1.600803348376281; H = 14.18296960958322 }
    6  = @{ D = 120512.9628395726;  E = -0.04447367095732655; F = 0.2617562671241508;  G = -1.145252808459844; H = 9.71280430611777 }
    8  = @{ D = 122164.7007662122;  E = -0.05677251663461622; F = 0.2245554879481218;  G = -0.9425775687075879; H = 7.102346691011166 }
    9  = @{ D = 123733.4508539603;  E = -0.1118279067509251;  F = 0.4943128361249983;  G = -2.68666777456219;  H = 15.14159420746112 }
    10 = @{ D = 125152.0335512803;  E = -0.1225329198767493;  F = 0.4432046915675618;  G = -1.929694107042016; H = 9.832829039523501 }
    11 = @{ D = 127217.7681993424;  E = -0.200059168140949;   F = 0.7849397096575603;  G = -2.642683839696853; H = 12.98399582730444 }
    13 = @{ D = 117650.839870485;   E = -0.03409780987832819; F = 0.1339026235787146;  G = -0.608512662289843; H = 6.747662536204578 }
    14 = @{ D = 117628.2773551512;  E = -0.0347753637637116;  F = 0.1420333826673406;  G = -0.9722330977683984; H = 9.295778709938382 }
    18 = @{ D = 117636.65421703;    E = -0.032172946883964;   F = 0.1381679263345998;  G = -1.357197505184592; H = 13.87796758682234 }
    20 = @{ D = 117705.7137032551;  E = -0.03076868531363822; F = 0.1530874067019525;  G = -0.5956323420837806; H = 6.325539329327718 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
